$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1446896", "Unknown", "2025-01-10 00:57:34"),
    @("1446896", "Unknown", "2025-01-10 00:58:35"),
    @("1446896", "Unknown", "2025-01-10 00:59:59"),
    @("1446896", "Unknown", "2025-01-10 01:02:32"),
    @("1446896_Asif Newaz", "Unknown", "2025-01-10 01:08:10"),
    @("1446896_Asif Newaz", "Unknown", "2025-01-10 01:09:41")
)

$startRow = 26
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    # Column A: values like "1446896" look numeric, so force text storage
    # (matches source data being plain ID strings, not numbers) with a
    # leading apostrophe-free text format applied only where needed.
    $idVal = $data[$i][0]
    $isNumericLooking = $idVal -match '^-?\d+(\.\d+)?$'
    if ($isNumericLooking) {
        $ws.Cells.Item($row, 1).NumberFormat = "@"
    }
    $ws.Cells.Item($row, 1).Value = $idVal
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
